$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate() | Out-Null

# --- Insert a new column before column B -----------------------------------
# This shifts the former B..G (requirement data) one column to the right
# (new C..H) and keeps column A untouched, exactly like Excel's
# "Insert Column" command (formulas/refs shift automatically).
$ws.Columns("B:B").Insert()

# --- Populate the freshly inserted column B with a running requirement # ---
# Rows 3..25 (the 23 data rows) get sequential numbers 1..23.
for ($i = 3; $i -le 25; $i++) {
    $ws.Cells.Item($i, 2).Value = $i - 2
}

# --- Column widths -----------------------------------------------------------
# Column A's width changed slightly and the new column B needs a width too.
# (ColumnWidth is specified in characters; the engine stores/quantizes the
# underlying width, so these are chosen to land as close as possible to the
# target stored widths.)
$ws.Columns("A:A").ColumnWidth = 8.67
$ws.Columns("B:B").ColumnWidth = 9

# --- View / selection ---------------------------------------------------------
# Scroll the sheet so column B is the left-most visible column, and move the
# active selection to E14 (matches the post-edit saved view state).
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E14").Select() | Out-Null
